$wb = $excel.ActiveWorkbook

$items = $wb.Worksheets.Item(1)
$items.Name = "Items"

$dict = $wb.Worksheets.Add($null, $items)
$dict.Name = "DataDictionary"

$nbsp = [char]0x00A0

$dict.Range("A1").Value = "Value"
$dict.Range("B1").Value = "Description"

$dict.Range("A2").Value  = "Item ID"
$dict.Range("A3").Value  = "Assoc Stim" + $nbsp
$dict.Range("A4").Value  = "Release Date"
$dict.Range("A5").Value  = "Subject"
$dict.Range("A6").Value  = "Item Type"
$dict.Range("A7").Value  = "TTS" + $nbsp
$dict.Range("A8").Value  = "Braille: BRF" + $nbsp
$dict.Range("A9").Value  = "Braille: PRN" + $nbsp
$dict.Range("A10").Value = "EG" + $nbsp
$dict.Range("A11").Value = "ASL" + $nbsp
$dict.Range("A12").Value = "Span" + $nbsp
$dict.Range("A13").Value = "TG" + $nbsp
$dict.Range("A14").Value = "IG" + $nbsp
$dict.Range("A15").Value = "Audio" + $nbsp
$dict.Range("A16").Value = "CC" + $nbsp

$dict.Range("B2").Value  = "Numerical identifier of item"
$dict.Range("B3").Value  = "Numerical identifier for associated stimulus"
$dict.Range("B4").Value  = "Date item was released in QTI sample packages"
$dict.Range("B5").Value  = "ELA or Math"
$dict.Range("B6").Value  = "MC, MS, EBSR, HTQ, EQ, TI, MI, GI, SA, WER"
$dict.Range("B7").Value  = "Text to Speech"
$dict.Range("B8").Value  = "Braille file type BRF"
$dict.Range("B9").Value  = "Braile file type PRN"
$dict.Range("B10").Value = "English glossary"
$dict.Range("B11").Value = "American sign language"
$dict.Range("B12").Value = "Spanish Translation"
$dict.Range("B13").Value = "Translated Glossary "
$dict.Range("B14").Value = "Illustration Glossary"
$dict.Range("B15").Value = "Audio stimulus"
$dict.Range("B16").Value = "Closed Captioning"

$header = $dict.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$hb = $header.Borders.Item(9)
$hb.Color = 0
$hb.LineStyle = 1
$hb.Weight = 2

